# Apply the edits described by the commit "Add files via upload":
#  - Student ids (column B) for rows 3-5 are replaced with new placeholder ids
#  - The e-mail/hyperlink in G5 (row 5) is removed (cell is cleared, its
#    hyperlink is deleted) while the hyperlinks on G3/G4 are kept
#  - The selected/active cell is moved from B16 to B12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the student id numbers in column B (rows 3-5) ---------------
$ws.Range("B3").Value = 1234567891
$ws.Range("B4").Value = 1234567892
$ws.Range("B5").Value = 1234567893

# --- Remove the Gmail hyperlink/value that used to sit in G5 ------------
# The worksheet-level Hyperlinks collection only supports bulk delete in
# this runtime, so drop every hyperlink and recreate the two that must
# stay (G3 and G4) pointing at their original mail addresses.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G3"), "mailto:tmot50473@gmail.com")
$ws.Hyperlinks.Add($ws.Range("G4"), "mailto:tranvankhanh4572@gmail.com")

# G5 itself becomes an empty cell (keeps its existing cell formatting).
$ws.Range("G5").ClearContents()

# --- Restore the previously selected cell --------------------------------
$ws.Range("B12").Select()
